# Apply fix: reset "Diferencia Stock" (column L) values to 0 for a set of rows,
# and update the "Total_Ajuste_Stock:" total in C96 to reflect the new sum (0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,4,6,8,10,19,24,26,27,28,29,33,39,45,46,48,49,50,52,57,61,66,74,77,78,80,81)

foreach ($r in $rows) {
    $ws.Range("L$r").Value = 0
}

# Update the total at the bottom of the sheet (row 96, column C) which sums column L.
$ws.Range("C96").Value = 0
